$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed symbol list
# snapshot taken on Mon Jan 16 07:39:18 UTC 2023.
# Values are prefixed with a leading apostrophe so Excel stores them as
# text (matching the original inline-string cell type) rather than
# re-typing the cell as a number/percentage.

$ws.Range("D2").Value = "'303.13"
$ws.Range("E2").Value = "'2.10%"
$ws.Range("D3").Value = "'31.73"
$ws.Range("E3").Value = "'0.38%"
$ws.Range("D4").Value = "'5.161"
$ws.Range("E4").Value = "'0.73%"
$ws.Range("D5").Value = "'0.07814"
$ws.Range("E5").Value = "'4.25%"
$ws.Range("D6").Value = "'2.311"
$ws.Range("E6").Value = "'34.82%"
$ws.Range("D7").Value = "'7.949"
$ws.Range("E7").Value = "'2.72%"
$ws.Range("D8").Value = "'3.867"
$ws.Range("E8").Value = "'1.75%"
$ws.Range("D9").Value = "'0.9077"
$ws.Range("E9").Value = "'-2.63%"
$ws.Range("D10").Value = "'0.1735"
$ws.Range("E10").Value = "'2.31%"
$ws.Range("D11").Value = "'0.07339"
$ws.Range("E11").Value = "'2.14%"
$ws.Range("D12").Value = "'0.08177"
$ws.Range("E12").Value = "'3.25%"
$ws.Range("D13").Value = "'0.03010"
$ws.Range("E13").Value = "'-0.41%"
$ws.Range("D14").Value = "'0.09948"
$ws.Range("E14").Value = "'0.52%"
$ws.Range("D15").Value = "'0.001515"
$ws.Range("E15").Value = "'0.24%"
$ws.Range("D16").Value = "'0.006070"
$ws.Range("E16").Value = "'-3.69%"
$ws.Range("E17").Value = "'1.35%"
$ws.Range("D18").Value = "'2.242"
$ws.Range("E18").Value = "'0.63%"
$ws.Range("E19").Value = "'-1.25%"
$ws.Range("D20").Value = "'0.1337"
$ws.Range("E20").Value = "'0.69%"
$ws.Range("D21").Value = "'4.670"
$ws.Range("E21").Value = "'2.51%"
$ws.Range("D22").Value = "'0.04654"
$ws.Range("E22").Value = "'0.19%"
$ws.Range("E23").Value = "'0.30%"
$ws.Range("D24").Value = "'0.001263"
$ws.Range("E24").Value = "'3.73%"
$ws.Range("D25").Value = "'0.004521"
$ws.Range("E25").Value = "'2.14%"
$ws.Range("E26").Value = "'3.73%"
$ws.Range("E27").Value = "'46.06%"
$ws.Range("D39").Value = "'0.01821"
$ws.Range("E39").Value = "'8.92%"
$ws.Range("D40").Value = "'0.04567"
$ws.Range("E40").Value = "'2.35%"
$ws.Range("D41").Value = "'0.007290"
$ws.Range("E41").Value = "'2.99%"
$ws.Range("D42").Value = "'0.1362"
$ws.Range("E42").Value = "'2.68%"
$ws.Range("D43").Value = "'0.002239"
$ws.Range("E43").Value = "'8.62%"
$ws.Range("D44").Value = "'0.01076"
$ws.Range("E44").Value = "'-5.34%"
$ws.Range("D45").Value = "'0.00006539"
$ws.Range("E45").Value = "'8.94%"
$ws.Range("E47").Value = "'-57.48%"
